# Applies the "Saldo.xlsx" update:
#  - EDNA's balance changes from 110519.39 to 22488
#  - ANTONIO (005000645) row removed
#  - AYRTON (001000882) row removed
#  - ANNA (004691225) row removed
#  - New row for BRASFORT (004352384 / 3.13) inserted just above the
#    004181486 / ANDREA row
#  - Trailing block of 4 rows (GILSON, ANA, NATALIA, BRASFORT w/ -40005.47)
#    removed from the end of the list
#
# Row numbers below refer to the ORIGINAL (before-edit) layout and are
# processed from the bottom of the sheet upward so that earlier deletes /
# inserts never invalidate the row numbers used by later steps.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# 1) Remove the trailing GILSON / ANA / NATALIA / BRASFORT(-40005.47) rows
#    (original rows 239-242).
$ws.Rows("239:242").Delete()

# 2) Insert a new row above the 004181486 / ANDREA row (original row 201)
#    and populate it with the BRASFORT / 3.13 account.
$ws.Rows.Item(201).Insert()
$ws.Range("A201").Value = "'004352384"
$ws.Range("B201").Value = "BRASFORT"
$ws.Range("C201").Value = 3.13

# 3) Remove the ANNA row (original row 85).
$ws.Rows.Item(85).Delete()

# 4) Remove the AYRTON row (original row 12).
$ws.Rows.Item(12).Delete()

# 5) Remove the ANTONIO row (original row 4).
$ws.Rows.Item(4).Delete()

# 6) Update EDNA's balance (row 3, column C) from 110519.39 to 22488.
$ws.Range("C3").Value = 22488
